$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new item "fuwen-aisi" is being inserted into the shop list as row 14.
# Everything that used to live in rows 14-37 (the Item name in column B,
# carrying its per-row style) logically slides down to rows 15-38, and one
# brand-new trailing row (38) is appended (continuing the Id sequence and
# reusing the last existing Item name).  The Id (A), Shelf (C, except the
# single boundary cell C22) and UseDiamond (D) columns do NOT shift - they
# stay put since they already happen to be sequential/constant per physical
# row - only column B (and the one shelf-boundary cell C22) changes.
#
# We shift column B by copying the B14:B37 block through a scratch column
# (to dodge self-overlap copy quirks) into B15:B38, keeping each cell's
# border style attached, then stamp the brand new B14 value in place.
# ---------------------------------------------------------------------------

$ws.Range("B14:B37").Copy($ws.Range("H14:H37"))
$ws.Range("H14:H37").Copy($ws.Range("B15:B38"))
$ws.Range("H14:H37").Clear()

# New row's item name (becomes a new shared string: "fuwen-aisi")
$ws.Range("B14").Value = "fuwen-aisi"

# Shelf boundary moved down by one row (old shelf-2/shelf-3 split was at
# row 22; it's now at row 23), so row 22 becomes shelf 2.
$ws.Range("C22").Value = 2

# Fill in the brand new trailing row 38 - Id continues the sequence,
# Shelf/UseDiamond copied from the row above (row 37) so the literal text
# ("true") keeps its shared-string type/style instead of being
# reinterpreted as a boolean or losing its border style.
$ws.Range("A38").Value = 15000045
$ws.Range("C37").Copy($ws.Range("C38"))
$ws.Range("C38").Value = 3
$ws.Range("D37").Copy($ws.Range("D38"))

# Table now spans one extra row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:D38"))

# Match the author's final cursor position from the diff.
$ws.Range("C14").Select()
